$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Ra"/"N" column pairs (V:W, X:Y, ...) were merged so the first
# cell of each pair carried the header/value and the second (right)
# cell of the pair was just an empty merged spacer. Column W was one
# such redundant/empty merged spacer cell - remove it outright so the
# columns after it (X, Y, Z) shift left by one (to W, X, Y) and the
# merges collapse away.
$ws.Range("W1").EntireColumn.Delete()

# Restore the selection to the cell that is now in the workbook's last
# used column (previously Z1, now W1 after the deletion above).
$ws.Range("W1").Select() | Out-Null

# Increase the indent of the "1_Ra" column header (C1) from 4 to 7.
$ws.Range("C1").IndentLevel = 7
